# Refresh the cryptos price/volume table (and the OKB/Bittensor rank swap)
# to match the latest scrape, as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.920.44"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "3.390.79"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.81"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.36"
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("D8").Value = "3.384.44"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +12.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.631"
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.62"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("E13").Value = "  +6.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.12"
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").Value = "3.926.66"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.32"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "3.389.69"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").Value = "64.811.01"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.66"
$ws.Range("E22").Value = "  +15.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("E23").Value = "  +13.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.14"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.10"
$ws.Range("E25").Value = "  +5.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.52"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.92"
$ws.Range("E27").Value = "  +7.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.81"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.78"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.03"
$ws.Range("E30").Value = "  +7.33%  "
$ws.Range("E31").Value = "  +4.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.53"
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "571.49"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.38"
$ws.Range("E34").Value = "  +6.40%  "
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.53"
$ws.Range("E38").Value = "  -4.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.67"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").Value = "0.0₃0753"
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.370"
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("D42").Value = "3.092.26"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("E45").Value = "  +4.41%  "
$ws.Range("E46").Value = "  +6.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.46"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.18"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.93"
$ws.Range("E50").Value = "  +3.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.34"
$ws.Range("E51").Value = "  +4.60%  "
